$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (SARWAR): strip currency symbols, bump the date ---
# Amount/Due Amount/Date are stored as plain text in the sheet, so a leading
# apostrophe keeps Excel from re-interpreting the numeric/date-looking text.
$ws.Range("B2").Value = "'2500"
$ws.Range("C2").Value = "'500"
$ws.Range("D2").Value = "'2026-02-02"

# --- Append 4 more dummy recipients (rows 3-6) ---
$data = @(
    @("John Doe",     "5000", "0",    "2026-02-01", "Service fee",     2002),
    @("Jane Smith",   "3500", "1500", "2026-01-30", "Partial payment", 2003),
    @("Ahmed Khan",   "4200", "4200", "2026-02-02", "Invoice #001",    2004),
    @("Maria Garcia", "1800", "0",    "2026-01-28", "Retainer",        2005)
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = "'" + $row[1]
    $ws.Cells.Item($r, 3).Value = "'" + $row[2]
    $ws.Cells.Item($r, 4).Value = "'" + $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}
